$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9298990368843079
$ws.Range("B1").Value = 1.314310073852539
$ws.Range("C1").Value = 3.579891681671143
$ws.Range("D1").Value = 3.042361974716187
$ws.Range("E1").Value = 0.564439058303833
